$d = $word.ActiveDocument

# --- Step 1: change the text of the last run of paragraph 2 -----------------
$p2 = $d.Paragraphs(2)
$p2.Range.Find.Execute("version de prueba 002", $true, $false, $false, $false, `
                        $false, $true, 1, $false, "entregable de prueba", 2)

# --- Step 2: insert a brand new paragraph right after paragraph 2 -----------
#     "1.2) Actividad 2" <br/> "1.24) Entregable de actividad 2"
$p2 = $d.Paragraphs(2)
$endOfP2 = $p2.Range.End
$ins = $d.Range($endOfP2, $endOfP2)
$ins.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$start3 = $d.Range($p3.Range.Start, $p3.Range.Start)
$start3.InsertAfter("1.2) Actividad 2")

$p3 = $d.Paragraphs(3)
$brPos = $p3.Range.End - 1
$brAndText = $d.Range($brPos, $brPos)
$breakChar = [string][char]11
$brAndText.InsertAfter($breakChar + "1.24) Entregable de actividad 2")

# Force the break+text portion into its own <w:r> instead of merging back
# with the preceding identically-formatted run.
$p3 = $d.Paragraphs(3)
$secondRun = $d.Range($brPos, $p3.Range.End - 1)
$secondRun.Font.Bold = 1
$secondRun.Font.Bold = 0

# that new paragraph keeps the tight (0/0) spacing seen in the target
$p3 = $d.Paragraphs(3)
$p3.SpaceBefore = 0
$p3.SpaceAfter = 0

# --- Step 3: insert another new paragraph after that one --------------------
#     <br/> "esto es lo que se necesita"  (back to the normal 240/240 spacing)
$p3 = $d.Paragraphs(3)
$endOfP3 = $p3.Range.End
$ins2 = $d.Range($endOfP3, $endOfP3)
$ins2.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$start4 = $d.Range($p4.Range.Start, $p4.Range.Start)
$start4.InsertAfter($breakChar + "esto es lo que se necesita")

$p4 = $d.Paragraphs(4)
$p4.SpaceBefore = 12
$p4.SpaceAfter = 12

Write-Output $d.Content.Text
